# Generate Report for Handback
# The "fc40310c-d3e2-430e-be7c-916464ffd6d6.md" file has now been handed
# back (in sync with en-US) for both the zh-cn and de-de locales. Update
# the Overview summary sheet plus the per-locale detail sheets to reflect
# the new status and record the handback timestamps.

$wb = $excel.ActiveWorkbook

$status = "Handed back: in sync with en-US"

# --- Overview sheet ---------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $status
$overview.Range("C3").Value = $status

# --- zh-cn detail sheet -------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $status
$zhcn.Range("H3").Value = "2016-03-24 21:02:12"

# --- de-de detail sheet -------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $status
$dede.Range("H3").Value = "2016-03-24 21:02:20"
